# Weekly price update: insert a new week's worth of data (2 rows) at the
# top of the "Zapallo italiano" price block for Agrícola del Norte S.A. de
# Arica, pushing the existing historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at row 526; this shifts the old rows 526-549 down
# to 528-551 and grows the used range to A1:R551 automatically.
$ws.Rows.Item(526).Insert()
$ws.Rows.Item(526).Insert()

# New row 526 - "Primera" quality, week of 2023-05-29 (serial 45075)
$ws.Cells.Item(526, 1).Value = 1
$ws.Cells.Item(526, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(526, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(526, 4).Value = 45075
$ws.Cells.Item(526, 5).Value = 15
$ws.Cells.Item(526, 6).Value = 100112032
$ws.Cells.Item(526, 7).Value = "Zapallo italiano"
$ws.Cells.Item(526, 8).Value = "Huracán"
$ws.Cells.Item(526, 9).Value = "Primera"
$ws.Cells.Item(526, 10).Value = 280
$ws.Cells.Item(526, 11).Value = 5000
$ws.Cells.Item(526, 12).Value = 6000
$ws.Cells.Item(526, 13).Value = 5643
$ws.Cells.Item(526, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(526, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(526, 16).Value = 81
$ws.Cells.Item(526, 17).Value = 70
$ws.Cells.Item(526, 18).Value = "Hortaliza"

# New row 527 - "Segunda" quality, same week
$ws.Cells.Item(527, 1).Value = 1
$ws.Cells.Item(527, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(527, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(527, 4).Value = 45075
$ws.Cells.Item(527, 5).Value = 15
$ws.Cells.Item(527, 6).Value = 100112032
$ws.Cells.Item(527, 7).Value = "Zapallo italiano"
$ws.Cells.Item(527, 8).Value = "Huracán"
$ws.Cells.Item(527, 9).Value = "Segunda"
$ws.Cells.Item(527, 10).Value = 200
$ws.Cells.Item(527, 11).Value = 4000
$ws.Cells.Item(527, 12).Value = 5000
$ws.Cells.Item(527, 13).Value = 4500
$ws.Cells.Item(527, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(527, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(527, 16).Value = 45
$ws.Cells.Item(527, 17).Value = 100
$ws.Cells.Item(527, 18).Value = "Hortaliza"

# Keep the date-formatted style (s="2") consistent for column D on the new
# rows, matching every other row in this block.
$ws.Cells.Item(526, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(527, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
